# Apply the "gamedata" edit:
#  - B2 total points 163 -> 179
#  - Row 6 (Rolls) roll-history strings change for Nick/Matt/Jasper/Nolan
#  - Selection moves to E6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total points for Nick (row 2) goes up from 163 to 179
$ws.Range("B2").Value = 179

# New roll histories for row 6 ("Rolls"): Nick, Matt, Jasper, Nolan
$ws.Range("B6").Value = "5/12/5/16/"
$ws.Range("C6").Value = "2/"
$ws.Range("D6").Value = "7/25/3/"
$ws.Range("E6").Value = "3/3/3/3/3/"

# Leave the selection on E6, matching the saved view state
[void]$ws.Range("E6").Select()
